$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '332.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.89%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.36%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.711'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.24%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08121'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.57%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.057'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.28%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.754'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.06%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.536'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.49%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.939'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.15%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9280'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.60%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1278'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.28%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1966'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.82%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.817'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '14.30%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09393'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.55%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.03723'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '8.16%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '9.21%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001307'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.18%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006247'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.70%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.380'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.18%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3515'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.44%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1417'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.18%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2610'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '6.60%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04420'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.28%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001259'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.42%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004419'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.59%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.38%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02912'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '16.98%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05516'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.25%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007830'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.43%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009892'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '10.23%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1426'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.62%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002091'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.94%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01109'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '5.75%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006762'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.10%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.13%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002279'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '26.64%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002995'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '3.48%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.13%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.13%'
